# Add "exceldelib" (Housni Achbouq's row became a REINSCRIPTION record)
# and "ReinscInsc" (new REINSCRIPTION/INSCRIPTION rows) to the "Emp Info" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Emp Info")

# --- Row 2: student 100 becomes a REINSCRIPTION record ---
$ws.Range("B2").Value = "J133341333"
$ws.Range("C2").Value = "Achbouq"
$ws.Range("D2").Value = "Housni"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "REINSCRIPTION"

# --- Row 3: student 101 becomes a REINSCRIPTION record ---
$ws.Range("B3").Value = "K133341333"
$ws.Range("C3").Value = "Benabbou2"
$ws.Range("D3").Value = "Oussama2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = "REINSCRIPTION"

# --- Row 4: student 102 becomes a REINSCRIPTION record ---
$ws.Range("B4").Value = "M133341333"
$ws.Range("C4").Value = "Lfelous"
$ws.Range("D4").Value = "Rim"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = "REINSCRIPTION"

# --- Row 5: new student 103, INSCRIPTION ---
$ws.Range("A5").Value = 103
$ws.Range("B5").Value = "A133341333"
$ws.Range("C5").Value = "Yamani"
$ws.Range("D5").Value = "Jamal"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "INSCRIPTION"

# --- Row 6: new student 104, INSCRIPTION ---
$ws.Range("A6").Value = 104
$ws.Range("B6").Value = "B133341333"
$ws.Range("C6").Value = "Ferdous"
$ws.Range("D6").Value = "Kamal"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "INSCRIPTION"

# --- Row 7: new student 105, INSCRIPTION ---
$ws.Range("A7").Value = 105
$ws.Range("B7").Value = "C133341333"
$ws.Range("C7").Value = "Touhami"
$ws.Range("D7").Value = "Badr"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "INSCRIPTION"

# Move the active selection to E2, like in the saved workbook.
[void]$ws.Range("E2").Select()
